$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be
# auto-recognized as numbers by Excel, so they stay text like the rest
# of the (inline string) column.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "39.912.95"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "2.210.56"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "288.74"
$ws.Range("E5").Value = "  -1.90%  "
$ws.Range("D6").Value = "87.05"
$ws.Range("E6").Value = "  +2.85%  "
$ws.Range("D7").Value = "0.512"
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").Value = "30.42"
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("D11").Value = "0.0776"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  +2.43%  "
$ws.Range("D13").Value = "6.45"
$ws.Range("E13").Value = "  +1.94%  "
$ws.Range("D14").Value = "2.549.61"
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("D15").Value = "13.93"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").Value = "2.216.21"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "0.726"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "39.847.82"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "11.58"
$ws.Range("E19").Value = "  +10.17%  "
$ws.Range("D20").Value = "0.0₃0881"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").Value = "5.78"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "65.40"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "235.07"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").Value = "22.45"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").Value = "153.99"
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("D31").Value = "31.62"
$ws.Range("E31").Value = "  -2.94%  "
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("E33").Value = "  +2.14%  "
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("D35").Value = "2.38"
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("D36").Value = "2.83"
$ws.Range("E36").Value = "  +6.44%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "15.69"
$ws.Range("E38").Value = "  -2.16%  "
$ws.Range("D39").Value = "0.0985"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("D40").Value = "1.70"
$ws.Range("E40").Value = "  +2.94%  "
$ws.Range("D41").Value = "3.84"
$ws.Range("E41").Value = "  +3.74%  "
$ws.Range("D42").Value = "2.091.13"
$ws.Range("E42").Value = "  +7.30%  "
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0267"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "9.95"
$ws.Range("E45").Value = "  +6.02%  "
$ws.Range("D46").Value = "17.45"
$ws.Range("E46").Value = "  +8.11%  "
$ws.Range("E47").Value = "  +2.59%  "
$ws.Range("D48").Value = "2.423.20"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "1.44"
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "88.47"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("D51").Value = "68.55"
$ws.Range("E51").Value = "  -3.07%  "
